# This script re-fetches/re-generates the class schedule rows so that
# rows end up shuffled into their new (generated) order.
# It reads the current values for the affected rows, then writes them
# back out in their new positions (a like a re-fetch from the backend
# re-ordered the rows for the same day/time block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return @(
        $ws.Cells.Item($row, 2).Value(),  # Time
        $ws.Cells.Item($row, 3).Value(),  # Unit
        $ws.Cells.Item($row, 4).Value(),  # Classroom
        $ws.Cells.Item($row, 5).Value(),  # Lecturer
        $ws.Cells.Item($row, 6).Value()   # Delivery Mode
    )
}

function Set-RowData($row, $data) {
    $ws.Cells.Item($row, 2).Value = $data[0]
    $ws.Cells.Item($row, 3).Value = $data[1]
    $ws.Cells.Item($row, 4).Value = $data[2]
    $ws.Cells.Item($row, 5).Value = $data[3]
    $ws.Cells.Item($row, 6).Value = $data[4]
}

# Capture original values of every row that gets reshuffled
$row2 = Get-RowData 2
$row3 = Get-RowData 3
$row4 = Get-RowData 4
$row5 = Get-RowData 5
$row7 = Get-RowData 7
$row8 = Get-RowData 8
$row9 = Get-RowData 9
$row10 = Get-RowData 10
$row11 = Get-RowData 11
$row15 = Get-RowData 15
$row16 = Get-RowData 16

# Row 2 <-> Row 4 swap
Set-RowData 2 $row4
Set-RowData 4 $row2

# Rows 3 -> 5 -> 7 -> 3 cyclic rotation
Set-RowData 3 $row7
Set-RowData 5 $row3
Set-RowData 7 $row5

# Row 8 <-> Row 9 swap
Set-RowData 8 $row9
Set-RowData 9 $row8

# Row 10 <-> Row 11 swap
Set-RowData 10 $row11
Set-RowData 11 $row10

# Row 15 <-> Row 16 swap
Set-RowData 15 $row16
Set-RowData 16 $row15
